$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# CAN bus connector changed back to JST GH series (from PH series)
# Row 31 holds: A=Ref Des, B=Part Description, C=Mfr Part Number, D=Datasheet,
#               E=Price, F=Quantity, G=Total Cost (formula)

$ws.Range("C31").Value = "SM03B-GHS-TB(LF)(SN)"
$ws.Range("B31").Value = "CONN HEADER GH SIDE 3POS 1.25MM"
$ws.Range("D31").Value = "http://www.jst-mfg.com/product/pdf/eng/eGH.pdf"

# New crimp connector price, which recalculates the formula-driven total cost
$ws.Range("E31").Value = 0.68

# Update the active selection left by the editor, as recorded in the sheet view
$ws.Range("K35").Select()
